$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column before column L (12th column). This shifts the
# existing columns L..T (and their data in every row) one position to
# the right, becoming M..U, and leaves a blank column L behind.
$ws.Columns.Item(12).Insert()

# Row 1 holds the "raw" lower-case column labels and is NOT supposed to
# shift with the rest of the table - the insert operation above moved
# the old L1:T1 text into M1:U1, so put the original labels back in
# L1:T1, and add one further (duplicated) label in the brand-new U1 cell.
$ws.Range("L1").Value = "P/l before tax"
$ws.Range("M1").Value = "Tax"
$ws.Range("N1").Value = "P/l after tax from ordinary activities"
$ws.Range("O1").Value = "Net profit/(loss) for the period"
$ws.Range("P1").Value = "Equity share capital"
$ws.Range("Q1").Value = "Basic eps"
$ws.Range("R1").Value = "Diluted eps"
$ws.Range("S1").Value = "Basic eps."
$ws.Range("T1").Value = "Diluted eps."
$ws.Range("U1").Value = "Diluted eps."

# Row 2 holds the "nice" Title-Case column labels; the insert already
# shifted M2:U2 correctly, we just need to label the new column.
$ws.Range("L2").Value = "Exceptional Items"
